$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the new "TODO" sheet after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Regles de gestion"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "TODO"

# --- Update the reworded explanation texts on the "Regles de gestion" sheet ---
$ws1.Range("B2").Value = "Si l'engagement est au statut 'SAISI'. Avec l'utilisateur en cours a le droit de Valider P et si l'utilisateur en cours n'est pas celui qui a effectué la saisie"
$ws1.Range("B3").Value = "Si l'engagement est au statut 'VALIDP'. Avecl'utilisateur en cours a le droit de Valider S et si l'utilisateur en cours n'est pas celui qui a effectué la validation P"
$ws1.Range("B4").Value = "Si l'engagement est au statut 'VALIDP' ou 'VALIDS'. Avec l'utilisateur en cours a le droit de Valider F et si l'utilisateur en cours n'est pas celui qui a effectué la validation S"

# --- Fill the new "TODO" sheet with its data ---
$ws2.Cells.Item(1, 1).Value = 1
$ws2.Cells.Item(1, 2).Value = "Loguer les connexions"

$ws2.Cells.Item(2, 1).Value = 2
$ws2.Cells.Item(2, 2).Value = "Vue pour avoir les états des engagements par lignes budgétaire"

$ws2.Cells.Item(3, 1).Value = 3
$ws2.Cells.Item(3, 2).Value = "Validation des formulaires"

$ws2.Cells.Item(4, 1).Value = 4
$ws2.Cells.Item(4, 2).Value = "Apurements"

$ws2.Cells.Item(5, 1).Value = 5
$ws2.Cells.Item(5, 2).Value = "Gérer les validations, Imputations, apurement des Réalisations directes"

$ws2.Cells.Item(6, 1).Value = 6
$ws2.Cells.Item(6, 2).Value = "Backup journalier de la base de données."

$ws2.Cells.Item(7, 1).Value = 7
$ws2.Cells.Item(7, 2).Value = "Ajout de fichiers lors de création, imputation, apurement engagement"

# --- Match the selections / active sheet of the final workbook ---
$ws1.Range("B3").Select()
$ws2.Range("B2").Select()
$ws2.Activate()
